$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing Excel to keep it as literal
# text (matches the source file, which stores every data cell as a string -
# the numeric-looking "Tag"/"Amount" columns included - never as a real
# number). Only touch NumberFormat (and reset the resulting style) when the
# text actually needs the nudge, so untouched-looking cells don't pick up a
# stray style index.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    # Plain numbers (e.g. "10904603", "8.65") and slash-dates that Excel's
    # auto-detect could parse as a real date/number (e.g. "10/02/2022" reads
    # as a valid M/D/Y date) need the NumberFormat nudge so they land as
    # literal text, same as every other cell in this sheet.
    $looksNumeric = $text -match '^-?[0-9]+(\.[0-9]+)?$'
    $looksDate = $text -match '^[0-9]{1,2}/[0-9]{1,2}/[0-9]{4}$'
    if ($looksNumeric -or $looksDate) {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

# The two newest "Feb 2022ST.pdf" statement rows (previously rows 5 and 6)
# move up to become rows 3 and 4; the two old "test.pdf" rows (2 and 3) are
# gone, and the former row 4 becomes the new row 2. Net effect: rows 5 and 6
# are removed and rows 2-4 are rewritten with the statement's data.
$ws.Rows("5:6").Delete()

Set-TextValue 2 1 "Feb 2022ST.pdf"
Set-TextValue 2 2 "14/01/2022"
Set-TextValue 2 3 "10904603"
Set-TextValue 2 4 "15:55:47"
Set-TextValue 2 5 "8.65"

Set-TextValue 3 1 "Feb 2022ST.pdf"
Set-TextValue 3 2 "29/01/2022"
Set-TextValue 3 3 "10904603"
Set-TextValue 3 4 "19:23:19"
Set-TextValue 3 5 "8.56"

Set-TextValue 4 1 "Feb 2022ST.pdf"
Set-TextValue 4 2 "10/02/2022"
Set-TextValue 4 3 "10904603"
Set-TextValue 4 4 "15:08:04"
Set-TextValue 4 5 "8.56"
